$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 273, shifting rows 273:317 down to 274:318
$ws.Rows.Item(273).Insert()

# Populate the newly inserted row 273 with the new record
$ws.Cells.Item(273, 1).Value = 7
$ws.Cells.Item(273, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(273, 3).Value = "Ñuble"
$ws.Cells.Item(273, 4).Value = (Get-Date -Year 2023 -Month 1 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(273, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(273, 5).Value = 16
$ws.Cells.Item(273, 6).Value = 100112006
$ws.Cells.Item(273, 7).Value = "Repollo"
$ws.Cells.Item(273, 8).Value = "Crespo record"
$ws.Cells.Item(273, 9).Value = "Primera"
$ws.Cells.Item(273, 10).Value = 20000
$ws.Cells.Item(273, 11).Value = 1200
$ws.Cells.Item(273, 12).Value = 1300
$ws.Cells.Item(273, 13).Value = 1250
$ws.Cells.Item(273, 14).Value = "$/unidad"
$ws.Cells.Item(273, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(273, 16).Value = 1250
$ws.Cells.Item(273, 17).Value = 1
$ws.Cells.Item(273, 18).Value = "Hortaliza"
